$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Datos"

# Format the working range as Text before entering data
$ws.Range("A1:F2").NumberFormat = "@"
$ws.Range("A3:A9").NumberFormat = "@"

# Enter the new header/data values in the same order they were
# originally typed in (controls shared-string insertion order)
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "password"
$ws.Range("D1").Value = "id_cuenta"
$ws.Range("E1").Value = "cuenta"
$ws.Range("B2").Value = "mat.rojasa@gmail.com"
$ws.Range("C2").Value = "@Rojas651"
$ws.Range("F1").Value = "empresa"
$ws.Range("D2").Value = "562267260251"
$ws.Range("E2").Value = "Movistar Hogar/Negocio Internet"
$ws.Range("F2").Value = "Internet"
$ws.Range("A1").Value = "TestCase"
